$d = $word.ActiveDocument

# 1. Insert a new "invisibleseparator" paragraph just before the horizontal-rule
#    paragraph (the empty paragraph that holds the w:pict rectangle). It is the
#    first paragraph with no explicit style ("Normal") that is otherwise empty.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Normal" -and $p.Range.Text -eq [char]13) {
        $hrIndex = $i
        break
    }
}
$hr = $d.Paragraphs.Item($hrIndex)
$hr.Range.InsertBefore(" " + [char]13)
$sep = $d.Paragraphs.Item($hrIndex)
$sep.Style = "invisibleseparator"

# 2. Remove the "_h2o_keep_element" marker paragraphs: Node Start, Head
#    Separator, Head End (around the resource header) and Node End (the
#    final paragraph). Walk bottom-up so earlier indices stay valid.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Head End" -or $styleName -eq "Head Separator" -or $styleName -eq "Node Start" -or $styleName -eq "Node End") {
        $p.Range.Delete()
    }
}

# 3. Remove the now-unused custom paragraph styles.
$stylesToDelete = @("HeadEnd", "HeadSeparator", "NodeEnd", "NodeStart")
foreach ($styleName in $stylesToDelete) {
    $style = $d.Styles.Item($styleName)
    $style.Delete()
}
